# Commit: "metadata and cellranger scripts added"
#
# The spreadsheet metadata table is corrected:
#   - columns 2 and 3 are swapped so the header reads Organ | Condition | Mice_id
#     (previously Organ | Mice_ID | Condition)
#   - the mouse identifiers are renamed from CF1_B/CF2_B/CF7_B/CG8_B/CG9_B
#     to the shorter F1/F2/F3/G2/G3
#   - the four trailing blank-but-formatted rows (7-10) are dropped, shrinking
#     the used range back down to A1:C6
#   - the stray "Times New Roman" font override on the id column is removed
#     (those cells fall back to the sheet's normal Calibri font)
#   - the bestFit custom width that used to sit on column D now belongs to
#     column C
#   - the saved selection moves to D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: Organ | Condition | Mice_id -------------------------------
$ws.Cells.Item(1, 1).Value = "Organ"
$ws.Cells.Item(1, 2).Value = "Condition"
$ws.Cells.Item(1, 3).Value = "Mice_id"

# --- data rows: Organ | Condition | Mice_id ---------------------------------
$data = @(
    @("Brain", "Flight", "F1"),
    @("Brain", "Flight", "F2"),
    @("Brain", "Flight", "F3"),
    @("Brain", "Ground", "G2"),
    @("Brain", "Ground", "G3")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# --- drop the leftover blank formatted rows below the table ----------------
$ws.Range("A7:C10").Clear()

# --- header is bold Calibri; the id column loses its Times New Roman look --
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("B2:B6").ClearFormats()

# --- the bestFit width moves from column D onto column C -------------------
$ws.Columns.Item(3).ColumnWidth = 17

# --- final saved selection ---------------------------------------------------
$ws.Range("D9").Select()
